# Weekly CompStat update (49th Precinct) - new crime data collected
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
# "Volume 32   Number  13" -> "...14"  (A8)
$a8 = $ws.Range("A8")
$a8.Value = ($a8.Text -replace '13$', '14')

# "Report Covering the Week  3/24/2025  Through  3/30/2025" -> updated week (C9)
$c9 = $ws.Range("C9")
$c9.Value = (($c9.Text -replace '3/24/2025','3/31/2025') -replace '3/30/2025','4/6/2025')

# --- Precinct crime-stat table updates (rows 15-28) ---
$ws.Range("C15").Value = 1
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("E15").Value = 0
$ws.Range("I15").Value = 7
$ws.Range("J15").Value = 12
$ws.Range("K15").Value = -41.666666666666
$ws.Range("L15").Value = 75
$ws.Range("M15").Value = 40
$ws.Range("N15").Value = 16.666666666666
$ws.Range("C16").Value = 5
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 20
$ws.Range("H16").Value = -15
$ws.Range("I16").Value = 67
$ws.Range("J16").Value = 74
$ws.Range("K16").Value = -9.459459459459
$ws.Range("L16").Value = 3.076923076923
$ws.Range("M16").Value = -14.102564102564
$ws.Range("N16").Value = -66.831683168316
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 35
$ws.Range("G17").Value = 29
$ws.Range("H17").Value = 20.689655172413
$ws.Range("I17").Value = 107
$ws.Range("J17").Value = 112
$ws.Range("K17").Value = -4.464285714285
$ws.Range("L17").Value = 7
$ws.Range("M17").Value = 59.701492537313
$ws.Range("N17").Value = 62.121212121212
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 66.666666666666
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = 14.285714285714
$ws.Range("I18").Value = 78
$ws.Range("J18").Value = 49
$ws.Range("K18").Value = 59.183673469387
$ws.Range("L18").Value = 25.806451612903
$ws.Range("M18").Value = -4.878048780487
$ws.Range("N18").Value = -82.068965517241
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = -33.333333333333
$ws.Range("F19").Value = 56
$ws.Range("G19").Value = 77
$ws.Range("H19").Value = -27.272727272727
$ws.Range("I19").Value = 197
$ws.Range("J19").Value = 265
$ws.Range("K19").Value = -25.660377358490
$ws.Range("L19").Value = 26.282051282051
$ws.Range("M19").Value = 143.20987654321
$ws.Range("N19").Value = 37.762237762237
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 11
$ws.Range("E20").Value = -27.272727272727
$ws.Range("F20").Value = 35
$ws.Range("G20").Value = 33
$ws.Range("H20").Value = 6.060606060606
$ws.Range("I20").Value = 132
$ws.Range("J20").Value = 124
$ws.Range("K20").Value = 6.451612903225
$ws.Range("L20").Value = 18.918918918918
$ws.Range("M20").Value = 116.393442622951
$ws.Range("N20").Value = -70.666666666666
$ws.Range("C21").Value = 37
$ws.Range("D21").Value = 44
$ws.Range("E21").Value = -15.909090909090
$ws.Range("G21").Value = 179
$ws.Range("H21").Value = -10.614525139664
$ws.Range("I21").Value = 590
$ws.Range("J21").Value = 637
$ws.Range("K21").Value = -7.378335949764
$ws.Range("L21").Value = 18.236472945891
$ws.Range("M21").Value = 56.914893617021
$ws.Range("N21").Value = -54.823889739663
$ws.Range("D22").Value = 2
$ws.Range("G22").Value = 3
$ws.Range("J22").Value = 8
$ws.Range("K22").Value = -50
$ws.Range("L22").Value = -42.857142857142
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = -33.333333333333
$ws.Range("F23").Value = 10
$ws.Range("H23").Value = 25
$ws.Range("I23").Value = 30
$ws.Range("J23").Value = 32
$ws.Range("K23").Value = -6.25
$ws.Range("L23").Value = -16.666666666666
$ws.Range("M23").Value = 66.666666666666
$ws.Range("C24").Value = 25
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = -7.407407407407
$ws.Range("F24").Value = 124
$ws.Range("G24").Value = 124
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 365
$ws.Range("J24").Value = 420
$ws.Range("K24").Value = -13.095238095238
$ws.Range("L24").Value = -7.594936708860
$ws.Range("M24").Value = 62.946428571428
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = -11.111111111111
$ws.Range("F25").Value = 39
$ws.Range("H25").Value = -27.777777777777
$ws.Range("I25").Value = 119
$ws.Range("J25").Value = 194
$ws.Range("K25").Value = -38.659793814433
$ws.Range("L25").Value = -30.409356725146
$ws.Range("C26").Value = 18
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 260
$ws.Range("F26").Value = 54
$ws.Range("G26").Value = 40
$ws.Range("H26").Value = 35
$ws.Range("I26").Value = 170
$ws.Range("J26").Value = 133
$ws.Range("K26").Value = 27.819548872180
$ws.Range("L26").Value = 18.055555555555
$ws.Range("M26").Value = 5.590062111801
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 1
$ws.Range("H27").Value = -83.333333333333
$ws.Range("I27").Value = 9
$ws.Range("J27").Value = 13
$ws.Range("K27").Value = -30.769230769230
$ws.Range("L27").Value = 0
$ws.Range("C28").Value = 2
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 6
$ws.Range("H28").Value = 50
$ws.Range("I28").Value = 20
$ws.Range("J28").Value = 16
$ws.Range("K28").Value = 25
$ws.Range("L28").Value = 33.333333333333
Write-Output "Update complete"
